$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("shows")

# Fill column A (Category) for the new "Google" block first (rows 42-52)
for ($r = 42; $r -le 52; $r++) {
    $ws.Cells.Item($r, 1).Value = "Google"
}

# Reproduce the exact original authoring order so that the shared-strings
# table ends up in the same sequence as the source workbook.
$ws.Cells.Item(43, 2).Value = "Volume Up"
$ws.Cells.Item(42, 2).Value = "Volume Down"
$ws.Cells.Item(44, 2).Value = "Play Metal"

$ws.Cells.Item(45, 2).Value = "Play Disturbed"
$ws.Cells.Item(45, 3).Value = "Okay Google. Play Disturbed on Pandora."

$ws.Cells.Item(46, 2).Value = "Next Song"
$ws.Cells.Item(46, 3).Value = "Okay Google. Next Song."

$ws.Cells.Item(47, 2).Value = "Song Info"
$ws.Cells.Item(47, 3).Value = "Okay Google. What's this song?"

$ws.Cells.Item(48, 2).Value = "Weather"
$ws.Cells.Item(48, 3).Value = "Okay Google. What's the weather forecast?"

$ws.Cells.Item(49, 2).Value = "Time"
$ws.Cells.Item(49, 3).Value = "Okay Google. What time is it?"

$ws.Cells.Item(50, 2).Value = "Off"
$ws.Cells.Item(50, 3).Value = "Okay Google. Off."

$ws.Cells.Item(42, 3).Value = "Okay Google. Lower Volume by Ten Percent."
$ws.Cells.Item(43, 3).Value = "Okay Google. Increase Volume by Ten Percent."
$ws.Cells.Item(44, 3).Value = "Okay Google. Play Heavy Metal on Pandora."

$ws.Cells.Item(51, 2).Value = "Jokes"
$ws.Cells.Item(51, 3).Value = "Okay Google. Tell me a joke."

$ws.Cells.Item(52, 2).Value = "Facts"
$ws.Cells.Item(52, 3).Value = "Okay Google. Give me a random fun fact."

# Update the view to match the saved window state
$ws.Application.ActiveWindow.ScrollRow = 33
$ws.Range("K37").Select()
